$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: the 2022-01-17 entry is no longer "in progress" -> values that
# used to be stored as text (e.g. "54446.0") become real numbers.
$ws.Range("B14").Value = 54446
$ws.Range("C14").Value = -74
$ws.Range("E14").Value = -2041

# --- Row 15: brand-new "未完成" (not yet finalized) entry for 2022-01-18.
# Its numeric-looking figures are still stored as plain text, same pattern
# row 14 used to follow before it was finalized above. Force text storage by
# pre-formatting the cells as Text ("@") before assigning the values, then
# drop back to the Normal style so no stray formatting is left behind.
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "2022-01-18"
$ws.Range("A15").Style = "Normal"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "54446.0"
$ws.Range("B15").Style = "Normal"

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "252.0"
$ws.Range("C15").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.46%"
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0"
$ws.Range("E15").Style = "Normal"
